$p = $ppt.ActivePresentation

# Delete slide 4 ("The Single Page App (SPA)") - all later slides shift up
$p.Slides.Item(4).Delete()
